$d = $word.ActiveDocument

# --- Step 1: strike-through + shrink the long sentence from "Nótese" through
# "...cargar un nivel)" (the closing paren, but NOT the trailing period). ---
$rng1 = $d.Content
$found1 = $rng1.Find.Execute(
    "Nótese que a priori no sabemos el tamaño del tablero y habrá que determinarlo a partir del propio archivo de nivel, en una pre-lectura de dicho archivo (se harán dos lecturas de archivo para cargar un nivel)",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found1) {
    $rng1.Font.StrikeThrough = $true
    $rng1.Font.Size = 8
}

# --- Step 2: the period + space right after it keeps the smaller size but
# is NOT struck through (it becomes its own run). ---
$rng2 = $d.Content
$found2 = $rng2.Find.Execute(
    ". Si el archivo de entrada",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found2) {
    # Narrow the range down to just ". " (2 chars) at the start of the match.
    $periodRange = $d.Range($rng2.Start, $rng2.Start + 2)
    $periodRange.Font.Size = 8
}
